# Generate Report for Handoff
# Rotates the three source-file rows on each sheet (the file that is now
# "Ready for handoff" moves to the bottom) and refreshes its handoff time.
#
# NOTE: hyperlink target URLs (per cell position) are left exactly as they
# were before the edit - only the cell text / hyperlink display text and a
# handful of status/datetime values change, matching the source diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"

$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/86a34335-664e-4588-96aa-6a7b9b630a1f.md", "", "", "ffff63458986-f214-4733-b71f-da0820b802e6.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/ffff63458986-f214-4733-b71f-da0820b802e6.md", "", "", "ffffff49eca66f-f1b1-40e6-a3d9-3acae5e81d2d.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/ffffff49eca66f-f1b1-40e6-a3d9-3acae5e81d2d.md", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "2016-03-09 05:23:42"
$wsZh.Range("G2").Value = "2016-03-09 05:24:36"
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "2016-03-09 05:23:42"
$wsZh.Range("G3").Value = "2016-03-09 05:24:36"
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "2016-03-09 05:28:35"
$wsZh.Range("G4").Value = "2016-03-09 05:28:05"
$wsZh.Range("H4").Value = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/86a34335-664e-4588-96aa-6a7b9b630a1f.md", "", "", "ffff63458986-f214-4733-b71f-da0820b802e6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ebfd34cd787bbc69342f9b492f99d683122428b1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.zh-cn.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/69a9a9546f3f41ce4d92771d4f9fca46684d991f/e2e/86a34335-664e-4588-96aa-6a7b9b630a1f.md", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dcfb4d98927e82f8d1f0b1e428841a37487aaef5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.zh-cn.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/ffff63458986-f214-4733-b71f-da0820b802e6.md", "", "", "ffffff49eca66f-f1b1-40e6-a3d9-3acae5e81d2d.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92ede07565bef86de935acba1ad55db5785cf8fc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73405e75ca977f966eb0e21b41525e3d75106621/e2e/bc727baf-ea08-49db-9362-f9f6c46b9863.md", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76074e5bdbbd85f107181cdbc7e1ae8ff530a17c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/ffffff49eca66f-f1b1-40e6-a3d9-3acae5e81d2d.md", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/92ede07565bef86de935acba1ad55db5785cf8fc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/73405e75ca977f966eb0e21b41525e3d75106621/e2e/bc727baf-ea08-49db-9362-f9f6c46b9863.md", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76074e5bdbbd85f107181cdbc7e1ae8ff530a17c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.zh-cn.xlf", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "2016-03-09 05:23:45"
$wsDe.Range("G2").Value = "2016-03-09 05:24:41"
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "2016-03-09 05:23:45"
$wsDe.Range("G3").Value = "2016-03-09 05:24:41"
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "2016-03-09 05:28:38"
$wsDe.Range("G4").Value = "2016-03-09 05:28:11"
$wsDe.Range("H4").Value = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/86a34335-664e-4588-96aa-6a7b9b630a1f.md", "", "", "ffff63458986-f214-4733-b71f-da0820b802e6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce63e5893a768f09c22e9f1ddb1e4cf2ecac35e1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.de-de.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b3ebffa3b59cd91a966533c8858714d234bb97b7/e2e/86a34335-664e-4588-96aa-6a7b9b630a1f.md", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bb8cd54b25e49c9f8724285678c26f6d647cfa19/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.de-de.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/ffff63458986-f214-4733-b71f-da0820b802e6.md", "", "", "ffffff49eca66f-f1b1-40e6-a3d9-3acae5e81d2d.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14d76271cb9a1ddfabe8cc5adc5db9b5139e65b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8c58297243df87d2e8f47254459d7c12b3061fd4/e2e/bc727baf-ea08-49db-9362-f9f6c46b9863.md", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/046169ae2be913c6c209f6dfc6f22efb4ef78fe9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf", "", "", "bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/e2e/ffffff49eca66f-f1b1-40e6-a3d9-3acae5e81d2d.md", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/14d76271cb9a1ddfabe8cc5adc5db9b5139e65b5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/8c58297243df87d2e8f47254459d7c12b3061fd4/e2e/bc727baf-ea08-49db-9362-f9f6c46b9863.md", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/046169ae2be913c6c209f6dfc6f22efb4ef78fe9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc727baf-ea08-49db-9362-f9f6c46b9863.6d7f2e327c1dab3f0d2c39f1bc1e44113bf764d0.de-de.xlf", "", "", "86a34335-664e-4588-96aa-6a7b9b630a1f.872ee9eeb0abb1911b48b36a5632d2b430edd915.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4534ce065ad383ffe7e18b79bad363a57ffc827d/.localization-config", "", "", ".localization-config")
